# Fix the "Class meetings" schedule table: each row's text had spilled
# across cell boundaries by one character. Re-distribute the text to the
# correct cells (and italicize the two "Spring Break" / "Final due" topic
# entries) without touching the row-number / unrelated second table.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# row -> (Date, Topic, Due, Meeting, TopicItalic)
$rows = @(
    @{ Row = 2;  Date = "Jan 27";  Topic = "Intro; Dual Coding / Cognitive Load";       Due = "-";            Meeting = "Online";     Italic = $false },
    @{ Row = 3;  Date = "Feb 3";   Topic = "Smart Enough Cities, ch 1-4";                Due = "-";            Meeting = "Harvey 104"; Italic = $false },
    @{ Row = 4;  Date = "Feb 10";  Topic = "Multimedia Principle & Contiguity";          Due = "Multimedia 1"; Meeting = "Online";     Italic = $false },
    @{ Row = 5;  Date = "Feb 17";  Topic = "Smart Enough Cities, ch 5-7";                Due = "-";            Meeting = "Harvey 104"; Italic = $false },
    @{ Row = 6;  Date = "Feb 24";  Topic = "Modality & Redundancy";                      Due = "Multimedia 2"; Meeting = "Online";     Italic = $false },
    @{ Row = 7;  Date = "Mar 2";   Topic = "Briefings";                                  Due = "Briefing";     Meeting = "Harvey 104"; Italic = $false },
    @{ Row = 8;  Date = "Mar 9";   Topic = "Coherence, Personalization, Segmenting";     Due = "Multimedia 3"; Meeting = "Online";     Italic = $false },
    @{ Row = 9;  Date = "Mar 16";  Topic = "Spring Break";                               Due = $null;          Meeting = $null;        Italic = $true  },
    @{ Row = 10; Date = "Mar 23";  Topic = "Pitches";                                    Due = "Pitch Harvey"; Meeting = "Harvey 104"; Italic = $false },
    @{ Row = 11; Date = "Mar 30";  Topic = "Workshops";                                  Due = "Wk 1 & 2";     Meeting = "Online";     Italic = $false },
    @{ Row = 12; Date = "Apr 6";   Topic = "Workshops";                                  Due = "Wk 3 & 4";     Meeting = "Harvey 104"; Italic = $false },
    @{ Row = 13; Date = "Apr 13";  Topic = "Workshops";                                  Due = "Wk 5 & 6";     Meeting = "Online";     Italic = $false },
    @{ Row = 14; Date = "Apr 20";  Topic = "Midpoint Critique";                          Due = "Project prototype"; Meeting = "Harvey 104"; Italic = $false },
    @{ Row = 15; Date = "Apr 27";  Topic = "Studio Session";                             Due = "-";            Meeting = "Online";     Italic = $false },
    @{ Row = 16; Date = "May 4";   Topic = "Studio Session";                             Due = "-";            Meeting = "Harvey 104"; Italic = $false },
    @{ Row = 17; Date = "May 11";  Topic = "User Testing";                               Due = "-";            Meeting = "Online";     Italic = $false },
    @{ Row = 18; Date = "May 18";  Topic = "Final due";                                  Due = "Final Project"; Meeting = "Harvey 104"; Italic = $true  }
)

foreach ($r in $rows) {
    # Column 2: Date
    $cell = $t.Cell($r.Row, 2)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $r.Date

    # Column 3: Topic
    $cell = $t.Cell($r.Row, 3)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $r.Topic
    if ($r.Italic) {
        # Re-fetch the (now exactly-sized) cell range, select it, and
        # toggle italics through the Selection (Range.Font.Italic alone
        # does not persist in this runtime).
        $cell.Range.Select()
        $word.Selection.Font.Italic = $true
    }

    # Column 4: Due (may be unchanged / untouched for the Spring Break row)
    if ($null -ne $r.Due) {
        $cell = $t.Cell($r.Row, 4)
        $rng = $cell.Range
        $rng.End = $rng.End - 1
        $rng.Text = $r.Due
    }

    # Column 5: Meeting (may be unchanged / untouched for the Spring Break row)
    if ($null -ne $r.Meeting) {
        $cell = $t.Cell($r.Row, 5)
        $rng = $cell.Range
        $rng.End = $rng.End - 1
        $rng.Text = $r.Meeting
    }
}
